# Reduce CBS concentration below solubility limit.
# Re-date the 2019-12-02 ITC run as 2019-12-03, shift the CBS/buffer
# titration series down by one concentration step (inserting a lower,
# 0.01 mM first point and a 0.02 mM "buffer into CAII" control),
# duplicate the CBS-into-CAII block for a second CAII 2 sample, and
# append the former "final cleaning/control" rows (now rows 20-22) plus
# a new CBS into CAII 2 fourth-concentration row (row 15), extending the
# used range from A1:K14 to A1:K22.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("plate")

# Row 2
$ws.Range("A2").Value = "20191203a1.itc"
$ws.Range("B2").Value = "initial cleaning water titration"
$ws.Range("C2").Value = "Plates Clean.setup"
$ws.Range("D2").Value = "water5inj.inj"
$ws.Range("E2").Value = "Control"
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = "Plate1, A1"
$ws.Range("I2").Value = "Plate1, B1"
$ws.Range("K2").Value = "Plate1, A1"

# Row 3
$ws.Range("A3").Value = "20191203a2.itc"
$ws.Range("B3").Value = "water into water 1"
$ws.Range("C3").Value = "Plates Clean.setup"
$ws.Range("D3").Value = "ChoderaWaterWater.inj"
$ws.Range("E3").Value = "Control"
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = "Plate1, C1"
$ws.Range("I3").Value = "Plate1, D1"
$ws.Range("K3").Value = "Plate1, C1"

# Row 4
$ws.Range("A4").Value = "20191203a3.itc"
$ws.Range("B4").Value = "buffer into buffer 1"
$ws.Range("C4").Value = "Plates Clean.setup"
$ws.Range("D4").Value = "ChoderaWaterWater.inj"
$ws.Range("E4").Value = "Control"
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = "Plate1, E1"
$ws.Range("I4").Value = "Plate1, F1"
$ws.Range("K4").Value = "Plate1, E1"

# Row 5
$ws.Range("A5").Value = "20191203a4.itc"
$ws.Range("B5").Value = "buffer into CAII 1"
$ws.Range("C5").Value = "Plates Clean.setup"
$ws.Range("D5").Value = "ChoderaWaterWater.inj"
$ws.Range("E5").Value = "Control"
$ws.Range("F5").Value = 0.01
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = "Plate1, G1"
$ws.Range("I5").Value = "Plate1, H1"
$ws.Range("K5").Value = "Plate1, G1"

# Row 6
$ws.Range("A6").Value = "20191203a5.itc"
$ws.Range("B6").Value = "buffer into CAII 1"
$ws.Range("C6").Value = "Plates Clean.setup"
$ws.Range("D6").Value = "ChoderaWaterWater.inj"
$ws.Range("E6").Value = "Control"
$ws.Range("F6").Value = 0.02
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = "Plate1, A2"
$ws.Range("I6").Value = "Plate1, B2"
$ws.Range("K6").Value = "Plate1, A2"

# Row 7
$ws.Range("A7").Value = "20191203a6.itc"
$ws.Range("B7").Value = "buffer into CAII 1"
$ws.Range("C7").Value = "Plates Clean.setup"
$ws.Range("D7").Value = "ChoderaWaterWater.inj"
$ws.Range("E7").Value = "Control"
$ws.Range("F7").Value = 0.04
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = "Plate1, C2"
$ws.Range("I7").Value = "Plate1, D2"
$ws.Range("K7").Value = "Plate1, C2"

# Row 8
$ws.Range("A8").Value = "20191203a7.itc"
$ws.Range("B8").Value = "CBS into CAII 1"
$ws.Range("C8").Value = "Plates Quick.setup"
$ws.Range("D8").Value = "ChoderaHSA.inj"
$ws.Range("E8").Value = "Onesite"
$ws.Range("F8").Value = 0.01
$ws.Range("G8").Value = 0.3311005429401354
$ws.Range("H8").Value = "Plate1, E2"
$ws.Range("I8").Value = "Plate1, F2"
$ws.Range("K8").Value = "Plate1, E2"

# Row 9
$ws.Range("A9").Value = "20191203a8.itc"
$ws.Range("B9").Value = "CBS into CAII 2"
$ws.Range("C9").Value = "Plates Quick.setup"
$ws.Range("D9").Value = "ChoderaHSA.inj"
$ws.Range("E9").Value = "Onesite"
$ws.Range("F9").Value = 0.01
$ws.Range("G9").Value = 0.3311005429401354
$ws.Range("H9").Value = "Plate1, G2"
$ws.Range("I9").Value = "Plate1, H2"
$ws.Range("K9").Value = "Plate1, G2"

# Row 10
$ws.Range("A10").Value = "20191203a9.itc"
$ws.Range("B10").Value = "CBS into CAII 1"
$ws.Range("C10").Value = "Plates Quick.setup"
$ws.Range("D10").Value = "ChoderaHSA.inj"
$ws.Range("E10").Value = "Onesite"
$ws.Range("F10").Value = 0.02
$ws.Range("G10").Value = 0.5231503635202742
$ws.Range("H10").Value = "Plate1, A3"
$ws.Range("I10").Value = "Plate1, B3"
$ws.Range("K10").Value = "Plate1, A3"

# Row 11
$ws.Range("A11").Value = "20191203a10.itc"
$ws.Range("B11").Value = "CBS into CAII 2"
$ws.Range("C11").Value = "Plates Quick.setup"
$ws.Range("D11").Value = "ChoderaHSA.inj"
$ws.Range("E11").Value = "Onesite"
$ws.Range("F11").Value = 0.02
$ws.Range("G11").Value = 0.5231503635202742
$ws.Range("H11").Value = "Plate1, C3"
$ws.Range("I11").Value = "Plate1, D3"
$ws.Range("K11").Value = "Plate1, C3"

# Row 12
$ws.Range("A12").Value = "20191203a11.itc"
$ws.Range("B12").Value = "CBS into CAII 1"
$ws.Range("C12").Value = "Plates Quick.setup"
$ws.Range("D12").Value = "ChoderaHSA.inj"
$ws.Range("E12").Value = "Onesite"
$ws.Range("F12").Value = 0.04
$ws.Range("G12").Value = 0.8575285224941928
$ws.Range("H12").Value = "Plate1, E3"
$ws.Range("I12").Value = "Plate1, F3"
$ws.Range("K12").Value = "Plate1, E3"

# Row 13
$ws.Range("A13").Value = "20191203a12.itc"
$ws.Range("B13").Value = "CBS into CAII 2"
$ws.Range("C13").Value = "Plates Quick.setup"
$ws.Range("D13").Value = "ChoderaHSA.inj"
$ws.Range("E13").Value = "Onesite"
$ws.Range("F13").Value = 0.04
$ws.Range("G13").Value = 0.8575285224941928
$ws.Range("H13").Value = "Plate1, G3"
$ws.Range("I13").Value = "Plate1, H3"
$ws.Range("K13").Value = "Plate1, G3"

# Row 14
$ws.Range("A14").Value = "20191203a13.itc"
$ws.Range("B14").Value = "CBS into CAII 1"
$ws.Range("C14").Value = "Plates Quick.setup"
$ws.Range("D14").Value = "ChoderaHSA.inj"
$ws.Range("E14").Value = "Onesite"
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0.3311005429401354
$ws.Range("H14").Value = "Plate1, A4"
$ws.Range("I14").Value = "Plate1, B4"
$ws.Range("K14").Value = "Plate1, A4"

# Row 15
$ws.Range("A15").Value = "20191203a14.itc"
$ws.Range("B15").Value = "CBS into CAII 2"
$ws.Range("C15").Value = "Plates Quick.setup"
$ws.Range("D15").Value = "ChoderaHSA.inj"
$ws.Range("E15").Value = "Onesite"
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0.3311005429401354
$ws.Range("H15").Value = "Plate1, C4"
$ws.Range("I15").Value = "Plate1, D4"
$ws.Range("K15").Value = "Plate1, C4"

# Row 16
$ws.Range("A16").Value = "20191203a15.itc"
$ws.Range("B16").Value = "CBS into CAII 1"
$ws.Range("C16").Value = "Plates Quick.setup"
$ws.Range("D16").Value = "ChoderaHSA.inj"
$ws.Range("E16").Value = "Onesite"
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0.5231503635202742
$ws.Range("H16").Value = "Plate1, E4"
$ws.Range("I16").Value = "Plate1, F4"
$ws.Range("K16").Value = "Plate1, E4"

# Row 17
$ws.Range("A17").Value = "20191203a16.itc"
$ws.Range("B17").Value = "CBS into CAII 2"
$ws.Range("C17").Value = "Plates Quick.setup"
$ws.Range("D17").Value = "ChoderaHSA.inj"
$ws.Range("E17").Value = "Onesite"
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0.5231503635202742
$ws.Range("H17").Value = "Plate1, G4"
$ws.Range("I17").Value = "Plate1, H4"
$ws.Range("K17").Value = "Plate1, G4"

# Row 18
$ws.Range("A18").Value = "20191203a17.itc"
$ws.Range("B18").Value = "CBS into CAII 1"
$ws.Range("C18").Value = "Plates Quick.setup"
$ws.Range("D18").Value = "ChoderaHSA.inj"
$ws.Range("E18").Value = "Onesite"
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0.8575285224941928
$ws.Range("H18").Value = "Plate1, A5"
$ws.Range("I18").Value = "Plate1, B5"
$ws.Range("K18").Value = "Plate1, A5"

# Row 19
$ws.Range("A19").Value = "20191203a18.itc"
$ws.Range("B19").Value = "CBS into CAII 2"
$ws.Range("C19").Value = "Plates Quick.setup"
$ws.Range("D19").Value = "ChoderaHSA.inj"
$ws.Range("E19").Value = "Onesite"
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0.8575285224941928
$ws.Range("H19").Value = "Plate1, C5"
$ws.Range("I19").Value = "Plate1, D5"
$ws.Range("K19").Value = "Plate1, C5"

# Row 20
$ws.Range("A20").Value = "20191203a19.itc"
$ws.Range("B20").Value = "final cleaning water titration"
$ws.Range("C20").Value = "Plates Clean.setup"
$ws.Range("D20").Value = "water5inj.inj"
$ws.Range("E20").Value = "Control"
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = "Plate1, E5"
$ws.Range("I20").Value = "Plate1, F5"
$ws.Range("K20").Value = "Plate1, E5"

# Row 21
$ws.Range("A21").Value = "20191203a20.itc"
$ws.Range("B21").Value = "final water into water test 1"
$ws.Range("C21").Value = "Plates Clean.setup"
$ws.Range("D21").Value = "ChoderaWaterWater.inj"
$ws.Range("E21").Value = "Control"
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = "Plate1, G5"
$ws.Range("I21").Value = "Plate1, H5"
$ws.Range("K21").Value = "Plate1, G5"

# Row 22
$ws.Range("A22").Value = "20191203a21.itc"
$ws.Range("B22").Value = "final water into water test 2"
$ws.Range("C22").Value = "Plates Clean.setup"
$ws.Range("D22").Value = "ChoderaWaterWater.inj"
$ws.Range("E22").Value = "Control"
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = "Plate1, A6"
$ws.Range("I22").Value = "Plate1, B6"
$ws.Range("K22").Value = "Plate1, A6"
